# jira to ado script
# The "Jira API Token" column (E) held a sensitive token value for every
# data row. Strip those values out so the secret no longer ships in the
# workbook; the now-unused shared string is dropped automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E5").ClearContents()
